# Add the new VARIABLE_CONTACT_PERIMETER variable row to the CONDUCTOR_files sheet.
# This documents a new user-defined auxiliary input file used to assign a
# variable contact perimeter between conductor components.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("CONDUCTOR_files")
$ws.Activate()

$ws.Range("A16").Value = "VARIABLE_CONTACT_PERIMETER"
$ws.Range("B16").Value = "-"
$ws.Range("C16").Value = "string"
$ws.Range("D16").Value = "external file for variable contact perimeter. Valid extension .xlsx"
$ws.Range("E16").Value = "none"

# Mirror the author's final on-screen state: new row selected (whole row).
$ws.Rows.Item(16).Select()
